$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 1024
$ws.Range("C3").Value = 1978
$ws.Range("D3").Value = 4031
$ws.Range("E3").Value = 9775
$ws.Range("F3").Value = 10400
$ws.Range("G3").Value = 11700
$ws.Range("B4").Value = 134.217728
$ws.Range("C4").Value = 258.998272
$ws.Range("D4").Value = 528.482304
$ws.Range("E4").Value = 1281.359872
$ws.Range("F4").Value = 1357.90592
$ws.Range("G4").Value = 1529.872384
$ws.Range("B5").Value = 969.6799999999999
$ws.Range("C5").Value = 992.95
$ws.Range("D5").Value = 921.92
$ws.Range("E5").Value = 768.33
$ws.Range("F5").Value = 1351.15
$ws.Range("G5").Value = 2341.43
$ws.Range("B6").Value = 1762
$ws.Range("C6").Value = 1762
$ws.Range("D6").Value = 1893
$ws.Range("E6").Value = 1860
$ws.Range("F6").Value = 3359
$ws.Range("G6").Value = 8979
$ws.Range("B7").Value = 1893
$ws.Range("C7").Value = 2180
$ws.Range("D7").Value = 2024
$ws.Range("E7").Value = 2040
$ws.Range("F7").Value = 4490
$ws.Range("G7").Value = 14484
$ws.Range("B12").Value = 4753
$ws.Range("C12").Value = 9004
$ws.Range("D12").Value = 17200
$ws.Range("E12").Value = 32400
$ws.Range("F12").Value = 58100
$ws.Range("G12").Value = 87800
$ws.Range("B13").Value = 19.5035136
$ws.Range("C13").Value = 36.9098752
$ws.Range("D13").Value = 70.35944959999999
$ws.Range("E13").Value = 132.120576
$ws.Range("F13").Value = 238.026752
$ws.Range("G13").Value = 359.661568
$ws.Range("B14").Value = 209.85605
$ws.Range("C14").Value = 221.69688
$ws.Range("D14").Value = 228.52552
$ws.Range("E14").Value = 243.13351
$ws.Range("F14").Value = 265.25533
$ws.Range("G14").Value = 320.86286
$ws.Range("B15").Value = 1417.216
$ws.Range("C15").Value = 1466.368
$ws.Range("D15").Value = 1449.984
$ws.Range("E15").Value = 1499.136
$ws.Range("F15").Value = 1564.672
$ws.Range("G15").Value = 1646.592
$ws.Range("B16").Value = 1531.904
$ws.Range("C16").Value = 1581.056
$ws.Range("D16").Value = 1597.44
$ws.Range("E16").Value = 1662.976
$ws.Range("G16").Value = 3850.24
$ws.Range("B21").Value = 8827
$ws.Range("C21").Value = 10900
$ws.Range("D21").Value = 13700
$ws.Range("E21").Value = 16700
$ws.Range("F21").Value = 16700
$ws.Range("G21").Value = 14400
$ws.Range("B22").Value = 1156.579328
$ws.Range("C22").Value = 1428.160512
$ws.Range("D22").Value = 1795.162112
$ws.Range("E22").Value = 2182.086656
$ws.Range("F22").Value = 2184.183808
$ws.Range("G22").Value = 1881.145344
$ws.Range("B23").Value = 45.52
$ws.Range("C23").Value = 87.26000000000001
$ws.Range("D23").Value = 115.79
$ws.Range("E23").Value = 100.04
$ws.Range("F23").Value = 165.46
$ws.Range("G23").Value = 1088.79
$ws.Range("B24").Value = 56
$ws.Range("C24").Value = 112
$ws.Range("D24").Value = 200
$ws.Range("E24").Value = 118
$ws.Range("F24").Value = 165
$ws.Range("G24").Value = 3851
$ws.Range("B25").Value = 62
$ws.Range("C25").Value = 326
$ws.Range("D25").Value = 1012
$ws.Range("E25").Value = 416
$ws.Range("F25").Value = 1860
$ws.Range("G25").Value = 9634
$ws.Range("B30").Value = 142000
$ws.Range("C30").Value = 211000
$ws.Range("D30").Value = 353000
$ws.Range("E30").Value = 532000
$ws.Range("F30").Value = 432000
$ws.Range("G30").Value = 203000
$ws.Range("B31").Value = 580.911104
$ws.Range("C31").Value = 866.123776
$ws.Range("D31").Value = 1447.03488
$ws.Range("E31").Value = 2177.892352
$ws.Range("F31").Value = 1767.899136
$ws.Range("G31").Value = 830.4721919999999
$ws.Range("B32").Value = ""
$ws.Range("F32").Value = 10.99
$ws.Range("G32").Value = 101.95
$ws.Range("B33").Value = 5.856
$ws.Range("C33").Value = 9
$ws.Range("D33").Value = 13
$ws.Range("E33").Value = 9.152000000000001
$ws.Range("F33").Value = 13
$ws.Range("G33").Value = 371
$ws.Range("B34").Value = 6.24
$ws.Range("C34").Value = 11
$ws.Range("D34").Value = 15
$ws.Range("E34").Value = 31.36
$ws.Range("F34").Value = 125
$ws.Range("G34").Value = 979
$ws.Range("B39").Value = 13300
$ws.Range("C39").Value = 18000
$ws.Range("D39").Value = 25800
$ws.Range("E39").Value = 32600
$ws.Range("F39").Value = 33400
$ws.Range("G39").Value = 35700
$ws.Range("B40").Value = 1742.733312
$ws.Range("C40").Value = 2485.12512
$ws.Range("D40").Value = 3376.41472
$ws.Range("E40").Value = 4278.19008
$ws.Range("F40").Value = 4373.610496
$ws.Range("G40").Value = 4673.503232
$ws.Range("B41").Value = 67.19
$ws.Range("C41").Value = 87.53
$ws.Range("D41").Value = 130.26
$ws.Range("E41").Value = 180.89
$ws.Range("F41").Value = 382.45
$ws.Range("G41").Value = 681.88
$ws.Range("B42").Value = 265
$ws.Range("C42").Value = 314
$ws.Range("D42").Value = 469
$ws.Range("E42").Value = 611
$ws.Range("F42").Value = 1369
$ws.Range("G42").Value = 2671
$ws.Range("B43").Value = 355
$ws.Range("C43").Value = 363
$ws.Range("D43").Value = 611
$ws.Range("E43").Value = 1221
$ws.Range("F43").Value = 2540
$ws.Range("G43").Value = 7046
$ws.Range("B48").Value = 426000
$ws.Range("C48").Value = 537000
$ws.Range("D48").Value = 819000
$ws.Range("E48").Value = 895000
$ws.Range("F48").Value = 986000
$ws.Range("G48").Value = 995000
$ws.Range("B49").Value = 1742.733312
$ws.Range("C49").Value = 2199.912448
$ws.Range("D49").Value = 3355.4432
$ws.Range("E49").Value = 3664.77312
$ws.Range("F49").Value = 4037.0176
$ws.Range("G49").Value = 4074.766336
$ws.Range("B50").Value = 2.06662
$ws.Range("C50").Value = 3.07137
$ws.Range("D50").Value = 4.45406
$ws.Range("E50").Value = 6.818750000000001
$ws.Range("F50").Value = 12.33315
$ws.Range("G50").Value = 24.86924
$ws.Range("B51").Value = 0.358
$ws.Range("C51").Value = 0.358
$ws.Range("D51").Value = 0.434
$ws.Range("E51").Value = 0.482
$ws.Range("F51").Value = 0.532
$ws.Range("G51").Value = 0.716
$ws.Range("B52").Value = 66.048
$ws.Range("C52").Value = 116.224
$ws.Range("D52").Value = 230.4
$ws.Range("E52").Value = 292.864
$ws.Range("F52").Value = 309.248
$ws.Range("G52").Value = 1073.152
$ws.Range("B57").Value = 5361
$ws.Range("C57").Value = 8291
$ws.Range("D57").Value = 11100
$ws.Range("E57").Value = 13200
$ws.Range("F57").Value = 12300
$ws.Range("G57").Value = 12100
$ws.Range("B58").Value = 702.54592
$ws.Range("C58").Value = 1086.324736
$ws.Range("D58").Value = 1455.423488
$ws.Range("E58").Value = 1729.101824
$ws.Range("F58").Value = 1606.418432
$ws.Range("G58").Value = 1592.786944
$ws.Range("B59").Value = 91.3
$ws.Range("C59").Value = 85.56999999999999
$ws.Range("D59").Value = 113.12
$ws.Range("E59").Value = 130.06
$ws.Range("F59").Value = 181.05
$ws.Range("G59").Value = 1017.82
$ws.Range("B60").Value = 68
$ws.Range("C60").Value = 86
$ws.Range("D60").Value = 172
$ws.Range("E60").Value = 262
$ws.Range("F60").Value = 182
$ws.Range("G60").Value = 4113
$ws.Range("B61").Value = 76
$ws.Range("C61").Value = 94
$ws.Range("D61").Value = 215
$ws.Range("E61").Value = 709
$ws.Range("F61").Value = 2212
$ws.Range("G61").Value = 8094
$ws.Range("B66").Value = 120000
$ws.Range("C66").Value = 163000
$ws.Range("D66").Value = 279000
$ws.Range("E66").Value = 340000
$ws.Range("F66").Value = 356000
$ws.Range("G66").Value = 215000
$ws.Range("B67").Value = 489.684992
$ws.Range("C67").Value = 665.84576
$ws.Range("D67").Value = 1145.044992
$ws.Range("E67").Value = 1394.60608
$ws.Range("F67").Value = 1458.569216
$ws.Range("G67").Value = 879.755264
$ws.Range("E68").Value = ""
$ws.Range("F68").Value = 11.86
$ws.Range("G68").Value = 99.77
$ws.Range("B69").Value = 5.344
$ws.Range("C69").Value = 8
$ws.Range("E69").Value = 8.896000000000001
$ws.Range("F69").Value = 13
$ws.Range("G69").Value = 302
$ws.Range("B70").Value = 5.792
$ws.Range("C70").Value = 9
$ws.Range("E70").Value = 23.424
$ws.Range("F70").Value = 110
$ws.Range("G70").Value = 865
